# 完成 /api/message/list， /api/message/id/{id} 接口
# Adds two new rows (23, 24) to the API prototype sheet documenting the
# "获取消息列表" (get message list) and "获取消息" (get message / mark-as-read)
# endpoints, mirroring the layout of the existing /api/tweet/* rows above.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 23: 获取消息列表  /api/message/list
# ---------------------------------------------------------------------
$ws.Range("A23").Value = "获取消息列表"
$ws.Range("B23").Value = "/api/message/list"
$ws.Range("C23").Value = "get"
$ws.Range("E23").Value = "int receiverId      //接收者Id  必须"
$ws.Range("F23").Value = "data:[{id: 消息Id,`n date:消息发送日期,`n status: 消息状态（未读或已读)`n content: 消息文本内容}]"

$ws.Range("A23").Style = $ws.Range("A22").Style
$ws.Range("B23").Style = $ws.Range("B22").Style
$ws.Range("C23").Style = $ws.Range("C22").Style
$ws.Range("E23").Style = $ws.Range("E22").Style
$ws.Range("F23").Style = $ws.Range("F21").Style

$ws.Rows.Item(23).RowHeight = 54

# ---------------------------------------------------------------------
# Row 24: 获取消息（获取之后会将状态设置为已读) /api/message/id/{id}
# ---------------------------------------------------------------------
$ws.Range("A24").Value = "获取消息`n（获取之后会将状态设置为已读)"
$ws.Range("B24").Value = "/api/message/id/{id}"
$ws.Range("C24").Value = "get"
$ws.Range("E24").Value = "int id    //消息Id 必须 放在路径上"
$ws.Range("F24").Value = "data:{id: 消息Id,`n date:消息发送日期,`n status: 消息状态（未读或已读)`n content: 消息文本内容}"

$ws.Range("A24").Style = $ws.Range("A19").Style
$ws.Range("B24").Style = $ws.Range("B22").Style
$ws.Range("C24").Style = $ws.Range("C22").Style
$ws.Range("E24").Style = $ws.Range("E21").Style
$ws.Range("F24").Style = $ws.Range("F21").Style

$ws.Rows.Item(24).RowHeight = 54

# Rich text run on A24: "获取消息\n（" stays default, the trailing
# "获取之后会将状态设置为已读)" warning note is small + red.
$noteChars = $ws.Range("A24").Characters(7, 14)
$noteChars.Font.Color = 255
$noteChars.Font.Size = 8

# Rich text run on E24: "放在路径上" at the end is bold to call out that
# the id is a path parameter.
$pathChars = $ws.Range("E24").Characters(21, 5)
$pathChars.Font.Bold = $true
$pathChars.Font.ColorIndex = -4105

# ---------------------------------------------------------------------
# Sheet view bookkeeping to match where the author left the cursor after
# typing in the new rows.
# ---------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 18
$ws.Range("E24").Select()
